$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "1.0000 at 0.00"
$ws.Range("C2").Value = "1.0000 at -120.00"
$ws.Range("D2").Value = "1.0000 at 120.00"

$ws.Range("C3").Value = "0.9907 at -120.19"
$ws.Range("D3").Value = "0.9982 at 120.02"

$ws.Range("C4").Value = "0.9890 at -120.26"
$ws.Range("D4").Value = "0.9962 at 120.07"

$ws.Range("B5").Value = "0.9969 at -0.07"
$ws.Range("C5").Value = "0.9980 at -120.05"
$ws.Range("D5").Value = "0.9973 at 119.99"

$ws.Range("B6").Value = "0.9724 at -0.77"
$ws.Range("C6").Value = "0.9788 at -120.55"
$ws.Range("D6").Value = "0.9782 at 119.50"

$ws.Range("B7").Value = "0.9641 at -3.02"
$ws.Range("C7").Value = "1.0154 at -121.01"
$ws.Range("D7").Value = "0.9686 at 119.75"

$ws.Range("B8").Value = "0.9606 at -3.07"
$ws.Range("D8").Value = "0.9705 at 119.92"

$ws.Range("D9").Value = "0.9724 at 120.05"

$ws.Range("B10").Value = "0.9574 at -3.28"
$ws.Range("C10").Value = "1.0177 at -121.19"
$ws.Range("D10").Value = "0.9667 at 119.76"

$ws.Range("B11").Value = "0.9641 at -3.02"
$ws.Range("C11").Value = "1.0154 at -121.01"
$ws.Range("D11").Value = "0.9686 at 119.75"

$ws.Range("B12").Value = "0.9552 at -2.99"
